$wb = $excel.ActiveWorkbook

# New values for column F (rows 2-9), to be applied to the sheets that
# contain this data table ("展览" and "全部类型").
$values = @(353, 95, 1535, 22, 51, 133, 55, 358)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 6).Value = $values[$i]
    }
}
